$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003208871385164791
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2195982.726253733
